$d = $word.ActiveDocument

function Find-ParagraphIndex($searchText) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        # Paragraph text can end with a paragraph mark (CR), a page-break
        # form-feed, manual line breaks, etc. - strip all trailing control
        # characters before comparing so lookups are stable regardless of
        # what trailing marks are currently present.
        $t = $p.Range.Text.TrimEnd([char]13, [char]12, [char]11, [char]10, [char]7)
        if ($t -eq $searchText) {
            return $i
        }
    }
    return -1
}

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# 1) Remove the stray _GoBack bookmark that currently sits right after
#    "Public static void main(String[] args) {"
# ---------------------------------------------------------------------------
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# ---------------------------------------------------------------------------
# 2) Merge the "Continued on Next Page" paragraph (right-justified, italic,
#    with an explicit page break run) into the following "Enemy Class"
#    paragraph. The resulting single paragraph keeps the "Enemy Class"
#    paragraph's own formatting (left ind 360, bold, no italics, no
#    right-justify) and gains the relocated _GoBack bookmark right before
#    its run.
# ---------------------------------------------------------------------------
$contIdx = Find-ParagraphIndex("Continued on Next Page")
$pCont = $d.Paragraphs.Item($contIdx)
$pEnemy = $d.Paragraphs.Item($contIdx + 1)
$mergeRange = $d.Range($pCont.Range.Start, $pEnemy.Range.End)
$mergeRange.InsertXML('<w:p ' + $wNs + '><w:pPr><w:ind w:left="360"/><w:rPr><w:rFonts w:ascii="Georgia" w:hAnsi="Georgia" w:cs="Courier New"/><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:rFonts w:ascii="Georgia" w:hAnsi="Georgia" w:cs="Courier New"/><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Enemy Class</w:t></w:r></w:p>')

# ---------------------------------------------------------------------------
# 3) The page that used to break just before "Enemy Class" now naturally
#    breaks before "Testing" instead, so move the lastRenderedPageBreak
#    marker there (inside the same run as the "Testing" text).
# ---------------------------------------------------------------------------
$testIdx = Find-ParagraphIndex("Testing")
$pTest = $d.Paragraphs.Item($testIdx)
$testRange = $d.Range($pTest.Range.Start, $pTest.Range.End)
$testRange.InsertXML('<w:p ' + $wNs + '><w:pPr><w:spacing w:after="240"/><w:ind w:left="360"/><w:rPr><w:rFonts w:ascii="Georgia" w:hAnsi="Georgia"/><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Georgia" w:hAnsi="Georgia"/><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:lastRenderedPageBreak/><w:t>Testing</w:t></w:r></w:p>')

# ---------------------------------------------------------------------------
# 4) Drop the whole "Conclusion Questions (CSP Only)" section through to the
#    end of the document (the CSP-only conclusion questions + the GUI
#    history discussion paragraphs), leaving the "...attribute again."
#    paragraph as the last one before the sectPr.
# ---------------------------------------------------------------------------
$conclIdx = Find-ParagraphIndex("Conclusion Questions (CSP Only)")
$pConcl = $d.Paragraphs.Item($conclIdx)
$pLast = $d.Paragraphs.Item($d.Paragraphs.Count)
$tailRange = $d.Range($pConcl.Range.Start, $pLast.Range.End)
$tailRange.Delete()
